$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for D and E columns (prices & percentages are text-formatted)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '23.032.82'
$ws.Range("E2").Value = '  -1.41%  '
$ws.Range("D3").Value = '1.598.68'
$ws.Range("E3").Value = '  -0.41%  '
$ws.Range("E4").Value = '  -0.65%  '
$ws.Range("D5").Value = '0.9996'
$ws.Range("E5").Value = '  -0.89%  '
$ws.Range("D6").Value = '301.76'
$ws.Range("E6").Value = '  -0.75%  '
$ws.Range("D7").Value = '0.3780'
$ws.Range("E7").Value = '  -0.48%  '
$ws.Range("D8").Value = '0.3649'
$ws.Range("E8").Value = '  -1.11%  '
$ws.Range("D9").Value = '50.12'
$ws.Range("E9").Value = '  +0.97%  '
$ws.Range("D10").Value = '1.256'
$ws.Range("E10").Value = '  -2.10%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '0.08145'
$ws.Range("E11").Value = '  -0.10%  '
$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D12").Value = '0.9991'
$ws.Range("E12").Value = '  -1.38%  '
$ws.Range("D13").Value = '22.46'
$ws.Range("E13").Value = '  -2.46%  '
$ws.Range("D14").Value = '6.594'
$ws.Range("E14").Value = '  -1.41%  '
$ws.Range("D15").Value = '7.384'
$ws.Range("E15").Value = '  -2.70%  '
$ws.Range("D16").Value = '0.00001246'
$ws.Range("E16").Value = '  -1.64%  '
$ws.Range("D17").Value = '1.599.13'
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("D18").Value = '91.73'
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("D19").Value = '0.06808'
$ws.Range("E19").Value = '  -0.87%  '
$ws.Range("D20").Value = '18.21'
$ws.Range("E20").Value = '  -2.07%  '
$ws.Range("D21").Value = '6.538'
$ws.Range("E21").Value = '  -1.95%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  -0.72%  '
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").Value = '13.02'
$ws.Range("E23").Value = '  -1.04%  '
$ws.Range("B24").Value = 'WrappedBTC'
$ws.Range("C24").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D24").Value = '23.063.72'
$ws.Range("E24").Value = '  -1.28%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '2.361'
$ws.Range("E25").Value = '  -1.63%  '
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = '2.809'
$ws.Range("E26").Value = '  -3.78%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '21.08'
$ws.Range("E27").Value = '  -1.12%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = '149.33'
$ws.Range("E28").Value = '  -1.08%  '
$ws.Range("B29").Value = 'HuobiToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D29").Value = '5.239'
$ws.Range("E29").Value = '  -1.66%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = '134.15'
$ws.Range("E30").Value = '  +0.49%  '
$ws.Range("B31").Value = 'WEMIXTOKEN'
$ws.Range("C31").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D31").Value = '2.369'
$ws.Range("E31").Value = '  -3.05%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '6.896'
$ws.Range("E32").Value = '  -8.60%  '
$ws.Range("B33").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C33").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D33").Value = '1.774.18'
$ws.Range("E33").Value = '  -0.81%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '0.9625'
$ws.Range("E34").Value = '  +0.75%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.07608'
$ws.Range("E35").Value = '  -1.27%  '
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").Value = '10.31'
$ws.Range("E36").Value = '  +0.38%  '
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '6.260'
$ws.Range("E37").Value = '  -1.60%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.02708'
$ws.Range("E38").Value = '  -3.44%  '
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").Value = '0.2530'
$ws.Range("E39").Value = '  -1.05%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '0.08837'
$ws.Range("E40").Value = '  -0.68%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '1.367'
$ws.Range("E41").Value = '  -1.79%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.7064'
$ws.Range("E42").Value = '  -1.69%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '12.40'
$ws.Range("E43").Value = '  -3.44%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '15.16'
$ws.Range("E44").Value = '  -4.55%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '0.6635'
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("B46").Value = 'Frax'
$ws.Range("C46").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D46").Value = '0.9986'
$ws.Range("E46").Value = '  -0.80%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '2.295'
$ws.Range("E47").Value = '  -1.46%  '
$ws.Range("D48").Value = '3.992'
$ws.Range("E48").Value = '  -0.38%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '131.52'
$ws.Range("E49").Value = '  -0.18%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.07932'
$ws.Range("E50").Value = '  -1.23%  '
$ws.Range("B51").Value = 'Flow'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D51").Value = '1.229'
$ws.Range("E51").Value = '  +4.63%  '

# Reset D:E cell style back to default (no explicit style) while keeping text type
$ws.Range("B2").Copy()
$ws.Range("D2:E51").PasteSpecial(-4122)
$excel.CutCopyMode = 0
